# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the
# per-locale handback status sheets (zh-cn and de-de).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 02:19:11"
$wsZhCn.Range("H2").Value = "2016-03-19 02:19:29"
$wsZhCn.Range("E4").Value = "2016-03-19 02:19:11"
$wsZhCn.Range("H4").Value = "2016-03-19 02:19:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 02:19:15"
$wsDeDe.Range("H2").Value = "2016-03-19 02:19:34"
$wsDeDe.Range("E4").Value = "2016-03-19 02:19:15"
$wsDeDe.Range("H4").Value = "2016-03-19 02:19:34"
